$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived values for rows 2-10 (columns G,H,I,J,M,N,O,P,Q,R,S,T)
# Row 2
$ws.Range("G2").Value = 13.32779766666667
$ws.Range("H2").Value = 39.983393
$ws.Range("I2").Value = 0.1697233513642653
$ws.Range("J2").Value = 0.1697233513642653
$ws.Range("M2").Value = 1.027065666666667
$ws.Range("N2").Value = 3.081197
$ws.Range("O2").Value = 0.1589549798736964
$ws.Range("P2").Value = 0.1589549798736964
$ws.Range("Q2").Value = 13.68852339571344
$ws.Range("R2").Value = 123.196710561421
$ws.Range("S2").Value = 0.0269783719002031
$ws.Range("T2").Value = 0.02697837190020309

# Row 3
$ws.Range("G3").Value = 13.32779766666667
$ws.Range("H3").Value = 39.983393
$ws.Range("I3").Value = 0.1697233513642653
$ws.Range("J3").Value = 0.1697233513642653
$ws.Range("O3").Value = 0.4795485327500095
$ws.Range("P3").Value = 0.4795485327500095
$ws.Range("Q3").Value = 41.29666975608111
$ws.Range("R3").Value = 371.67002780473
$ws.Range("S3").Value = 0.08139058412014774
$ws.Range("T3").Value = 0.08139058412014774

# Row 4
$ws.Range("G4").Value = 13.32779766666667
$ws.Range("H4").Value = 39.983393
$ws.Range("I4").Value = 0.1697233513642653
$ws.Range("J4").Value = 0.1697233513642653
$ws.Range("M4").Value = 2.335759666666667
$ws.Range("N4").Value = 7.007279
$ws.Range("O4").Value = 0.3614964873762942
$ws.Range("P4").Value = 0.3614964873762941
$ws.Range("Q4").Value = 31.13053223529411
$ws.Range("R4").Value = 280.174790117647
$ws.Range("S4").Value = 0.06135439534391448
$ws.Range("T4").Value = 0.06135439534391447

# Row 5
$ws.Range("I5").Value = 0.5514955210569645
$ws.Range("J5").Value = 0.5514955210569645
$ws.Range("M5").Value = 1.027065666666667
$ws.Range("N5").Value = 3.081197
$ws.Range("O5").Value = 0.1589549798736964
$ws.Range("P5").Value = 0.1589549798736964
$ws.Range("Q5").Value = 44.47920266679866
$ws.Range("R5").Value = 400.3128240011879
$ws.Range("S5").Value = 0.08766295945004351
$ws.Range("T5").Value = 0.0876629594500435

# Row 6
$ws.Range("I6").Value = 0.5514955210569645
$ws.Range("J6").Value = 0.5514955210569645
$ws.Range("O6").Value = 0.4795485327500095
$ws.Range("P6").Value = 0.4795485327500095
$ws.Range("S6").Value = 0.2644688679410693
$ws.Range("T6").Value = 0.2644688679410693

# Row 7
$ws.Range("I7").Value = 0.5514955210569645
$ws.Range("J7").Value = 0.5514955210569645
$ws.Range("M7").Value = 2.335759666666667
$ws.Range("N7").Value = 7.007279
$ws.Range("O7").Value = 0.3614964873762942
$ws.Range("P7").Value = 0.3614964873762941
$ws.Range("Q7").Value = 101.1549027159907
$ws.Range("R7").Value = 910.394124443916
$ws.Range("S7").Value = 0.1993636936658518
$ws.Range("T7").Value = 0.1993636936658517

# Row 8
$ws.Range("H8").Value = 65.67520200000001
$ws.Range("I8").Value = 0.27878112757877
$ws.Range("J8").Value = 0.27878112757877
$ws.Range("M8").Value = 1.027065666666667
$ws.Range("N8").Value = 3.081197
$ws.Range("O8").Value = 0.1589549798736964
$ws.Range("P8").Value = 0.1589549798736964
$ws.Range("Q8").Value = 22.48424837519934
$ws.Range("R8").Value = 202.358235376794
$ws.Range("S8").Value = 0.04431364852344979
$ws.Range("T8").Value = 0.04431364852344978

# Row 9
$ws.Range("H9").Value = 65.67520200000001
$ws.Range("I9").Value = 0.27878112757877
$ws.Range("J9").Value = 0.27878112757877
$ws.Range("O9").Value = 0.4795485327500095
$ws.Range("P9").Value = 0.4795485327500095
$ws.Range("Q9").Value = 67.83234049591334
$ws.Range("R9").Value = 610.4910644632201
$ws.Range("S9").Value = 0.1336890806887924
$ws.Range("T9").Value = 0.1336890806887924

# Row 10
$ws.Range("H10").Value = 65.67520200000001
$ws.Range("I10").Value = 0.27878112757877
$ws.Range("J10").Value = 0.27878112757877
$ws.Range("M10").Value = 2.335759666666667
$ws.Range("N10").Value = 7.007279
$ws.Range("O10").Value = 0.3614964873762942
$ws.Range("P10").Value = 0.3614964873762941
$ws.Range("Q10").Value = 51.13382931059534
$ws.Range("R10").Value = 460.2044637953581
$ws.Range("S10").Value = 0.1007783983665279
$ws.Range("T10").Value = 0.1007783983665279
